$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11 (columns D, L, M, N, O, P, Q, R, S, T)
# This reorders the existing weekly price records (same underlying data,
# different row order) as part of the "Fruta / hortaliza, semanal" update.

$data = @(
    @{ Row=2;  D=44208; L="Especial";     M=70;  N=24000; O=24000; P=24000; Q="$/caja 15 kilos granel";     R="Región de O'Higgins"; S=1600; T=15 },
    @{ Row=3;  D=44411; L="Primera";      M=210; N=8000;  O=8000;  P=8000;  Q="$/bandeja 8 kilos";           R="Región de O'Higgins"; S=1000; T=8  },
    @{ Row=4;  D=44217; L="Primera";      M=55;  N=18000; O=18000; P=18000; Q="$/caja 18 kilos granel";     R="Región de O'Higgins"; S=1000; T=18 },
    @{ Row=5;  D=44511; L="Primera";      M=15;  N=22000; O=22000; P=22000; Q="$/caja 15 kilos granel";     R="Región de O'Higgins"; S=1467; T=15 },
    @{ Row=6;  D=44601; L="Primera";      M=30;  N=28000; O=28000; P=28000; Q="$/caja 18 kilos granel";     R="Región de O'Higgins"; S=1556; T=18 },
    @{ Row=7;  D=44427; L="Primera";      M=55;  N=7000;  O=7000;  P=7000;  Q="$/caja 15 kilos granel";     R="Región de O'Higgins"; S=467;  T=15 },
    @{ Row=8;  D=44392; L="Especial";     M=500; N=7000;  O=7000;  P=7000;  Q="$/bandeja 8 kilos";           R="Región de O'Higgins"; S=875;  T=8  },
    @{ Row=9;  D=44264; L="Calibre 100";  M=50;  N=20000; O=20000; P=20000; Q="$/caja 18 kilos embalada";   R="Región de O'Higgins"; S=1111; T=18 },
    @{ Row=10; D=44495; L="Primera";      M=50;  N=24000; O=24000; P=24000; Q="$/bandeja 10 kilos";          R="China";               S=2400; T=10 },
    @{ Row=11; D=44418; L="Especial";     M=100; N=8000;  O=8000;  P=8000;  Q="$/caja 15 kilos granel";     R="Región de O'Higgins"; S=533;  T=15 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value  = $item.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value = $item.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $item.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $item.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $item.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $item.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $item.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $item.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $item.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $item.T   # T: Kg / unidad
}
